# Add a new translated-line entry ("Good luck to you!" / "Удачи вам!")
# to the bottom of the table, turning the former last row (row 12, part
# of the um1606 entry) into a non-final row with a simple bottom border
# and appending a brand-new row 13 for um2407 with the "new group" top
# border styling (matching the look of rows 2/6/8/11 etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 currently uses the "first row of group" style (s=4/5, no top
# border). Since it will no longer be the last row of the table, give it
# the "continuation row" style (s=6/7, thin bottom border only) -- the
# same look already used by rows 4, 7 and 9. Row 9 is a perfect template.
$ws.Range("A9:E9").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)   # xlPasteFormats

# Row 13 is a brand-new "first row of group" entry. Row 11 already has
# exactly that styling (s=4/5), so copy its formatting as the template.
$ws.Range("A11:E11").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(13).RowHeight = 43.2

# Fill in row 13's values. Order matters for shared-string allocation:
# set the English line first, then the filename, then the Russian lines,
# so new shared strings land in the same order as the source workbook.
$ws.Cells.Item(13, 3).Value = " Good luck to you!"
$ws.Cells.Item(13, 1).Value = "SCRIPT/G01P03A/um2407.ssb"
$ws.Cells.Item(13, 4).Value = " Удачи вам!"
$ws.Cells.Item(13, 5).Value = " Ôäàœé âàí!"
$ws.Cells.Item(13, 2).Value = 119
